# This commit ("Vygenerovany file ve slozce: ...") only touches the
# <w:nsid w:val="..."/> GUID-like identifiers that Word stamps on each
# <w:abstractNum> entry in word/numbering.xml:
#
#   abstractNumId 990   : b44b2ced -> 8838391d
#   abstractNumId 991   : 6a7f9daa -> 2ea6ad3d
#   abstractNumId 99721 : 4ffc813a -> a4346194
#   abstractNumId 99722 : a7070901 -> 780a1550
#
# Nothing else in the package changes: same abstractNumId/numId values,
# same levels, same paragraphs/text/formatting. <w:nsid> is an internal,
# random "list signature" Word mints for its own bookkeeping - it is not
# surfaced anywhere on the Word object model (ListTemplate / ListFormat /
# List only expose Name, OutlineNumbered, ListLevels, ListID-of-the-numId,
# etc. - verified here via Get-Member against every list-related object),
# so it cannot be targeted with Find/Replace (it never appears in any
# Range's text) and there is no settable property for it. We still try
# every plausible property name defensively (in case a given host exposes
# more than this one does) but guard every attempt so a missing member
# never throws and never leaves the document in a half-edited state.

$d = $word.ActiveDocument

$map = @{
    990   = "8838391d"
    991   = "2ea6ad3d"
    99721 = "a4346194"
    99722 = "780a1550"
}

# abstractNum declaration order in word/numbering.xml for this document:
# 1 -> abstractNumId 0 (unchanged), 2 -> 990, 3 -> 991, 4 -> 99721, 5 -> 99722
$order = @(0, 990, 991, 99721, 99722)

function Try-SetNsid($obj, $value) {
    if ($null -eq $obj) { return }
    # Dynamic "$obj.$prop = $value" member access isn't supported by this
    # interpreter, so each plausible property name is attempted literally;
    # every attempt is independently guarded so an unimplemented member
    # never aborts the script.
    try { $obj.Nsid = $value } catch { }
    try { $obj.NSID = $value } catch { }
    try { $obj.nsid = $value } catch { }
    try { $obj.ListId = $value } catch { }
    try { $obj.ListID = $value } catch { }
    try { $obj.Id = $value } catch { }
    try { $obj.Tag = $value } catch { }
    try { $obj.Guid = $value } catch { }
    try { $obj.GUID = $value } catch { }
}

if ($d.ListTemplates.Count -ge $order.Count) {
    for ($i = 1; $i -le $order.Count; $i++) {
        $abstractNumId = $order[$i - 1]
        if ($map.ContainsKey($abstractNumId)) {
            $newValue = $map[$abstractNumId]
            try {
                $lt = $d.ListTemplates.Item($i)
                Try-SetNsid $lt $newValue
            } catch {
                # Collection/index not available - nothing more to do for
                # this abstractNum via the object model.
            }
        }
    }
}

# word/document.xml <w:numId> -> <w:abstractNumId> (from numbering.xml):
#   numId 1002 -> abstractNumId 99721
#   numId 1003 -> abstractNumId 99722
$numIdToAbstract = @{
    1002 = 99721
    1003 = 99722
}

# Also sweep every List actually bound to a paragraph, in case the host
# resolves "the list" to a richer object than the plain ListTemplate view.
foreach ($p in $d.Paragraphs) {
    try {
        $lf = $p.Range.ListFormat
        if ($lf.ListType -ne 0) {
            $numId = $lf.List.ListID
            if ($numIdToAbstract.ContainsKey($numId)) {
                $abstractNumId = $numIdToAbstract[$numId]
                if ($map.ContainsKey($abstractNumId)) {
                    $newValue = $map[$abstractNumId]
                    Try-SetNsid $lf.ListTemplate $newValue
                    Try-SetNsid $lf.List $newValue
                }
            }
        }
    } catch {
    }
}

Write-Output "nsid update attempted for abstractNum 990, 991, 99721, 99722"
